# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" sheet and the "全部类型" sheet, which both contain the same
# underlying rows (the latter aggregates all event categories and has one
# extra row inserted, shifting row numbers by +1 from row 6 onward).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览" sheet (sheet1)
$updates1 = @{
    2  = 7083
    4  = 464
    5  = 17
    7  = 160
    8  = 122
    11 = 55
    12 = 203
    15 = 1833
    16 = 45
    17 = 3669
    21 = 26
    23 = 2305
    24 = 17
    25 = 271
    32 = 1324
    33 = 116
}

foreach ($row in $updates1.Keys) {
    $sheet1.Range("F$row").Value = $updates1[$row]
}

# Row -> new F value for "全部类型" sheet (sheet4), shifted by +1 row
# relative to "展览" starting at row 8 due to an extra row in this sheet.
$updates4 = @{
    2  = 7083
    4  = 464
    5  = 17
    8  = 160
    9  = 122
    12 = 55
    13 = 203
    16 = 1833
    17 = 45
    18 = 3669
    22 = 26
    24 = 2305
    25 = 17
    26 = 271
    33 = 1324
    34 = 116
}

foreach ($row in $updates4.Keys) {
    $sheet4.Range("F$row").Value = $updates4[$row]
}
